# Auto-generated edit script: updates crypto price/volume table cells
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.125.86'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.858.32'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '471.23'
$ws.Range('E5').Value = '  +10.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.15'
$ws.Range('E6').Value = '  +10.20%  '
$ws.Range('E7').Value = '  +3.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.745'
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('E10').Value = '  -3.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000311'
$ws.Range('E11').Value = '  -9.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.47'
$ws.Range('E12').Value = '  +3.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.43'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.495.22'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.85'
$ws.Range('E15').Value = '  -6.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.894.58'
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.04'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('E19').Value = '  +5.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.320.82'
$ws.Range('E20').Value = '  +0.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '434.91'
$ws.Range('E21').Value = '  +4.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.89'
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.29'
$ws.Range('E23').Value = '  +6.02%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.02'
$ws.Range('E24').Value = '  +4.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.60'
$ws.Range('E25').Value = '  +10.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.87'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.15'
$ws.Range('E27').Value = '  +8.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.97'
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('E29').Value = '  +2.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '732.43'
$ws.Range('E30').Value = '  +1.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.84'
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('E32').Value = '  +6.48%  '
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '44.72'
$ws.Range('E34').Value = '  +13.90%  '
$ws.Range('E35').Value = '  +7.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '57.95'
$ws.Range('E36').Value = '  +3.90%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.53'
$ws.Range('E38').Value = '  -4.48%  '
$ws.Range('E39').Value = '  +3.57%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₃0687'
$ws.Range('E42').Value = '  -8.77%  '
$ws.Range('E43').Value = '  +3.30%  '
$ws.Range('E44').Value = '  +9.21%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.45'
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('E48').Value = '  +5.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.15'
$ws.Range('E49').Value = '  +4.57%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.89'
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.21'
$ws.Range('E51').Value = '  +1.05%  '
